$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 ("I0") and J1 ("IF") ---
# Set the text first, then clone H1's formatting (bold, centered, bordered
# header style) onto the two new header cells via copy/paste-special so the
# workbook reuses the existing style index instead of creating a new one.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

# --- New data columns I ("I0") and J ("IF") for rows 2-35 ---
$iVals = @(1,9,9,7,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$jVals = @(1,9,9,7,2,4,4,4,5,6,5,5,6,1,5,4,5,5,7,5,7,8,5,6,6,5,5,6,6,5,5,4,4,3)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
